# Weekly update: a new daily price record is inserted at row 57 (pushing the
# existing rows 57-88 down to 58-89), and the newly-opened row 57 is filled
# with the new day's data (same market/category/etc. as the rest of the
# sheet, but a new date and new volume/price figures).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 57..88 down to 58..89, leaving a blank row 57 in their place.
$ws.Rows.Item(57).Insert()

# Populate the new row 57 with the new record.
$ws.Range("A57").Value = 8
$ws.Range("B57").Value = "Terminal La Palmera de La Serena"
$ws.Range("C57").Value = "Coquimbo"
$ws.Range("D57").Value = 44489
$ws.Range("E57").Value = 4
$ws.Range("F57").Value = 100112001
$ws.Range("G57").Value = "Berenjena"
$ws.Range("H57").Value = "Sin especificar"
$ws.Range("I57").Value = "Primera"
$ws.Range("J57").Value = 500
$ws.Range("K57").Value = 8000
$ws.Range("L57").Value = 9000
$ws.Range("M57").Value = 8500
$ws.Range("N57").Value = "$/caja 60 unidades"
$ws.Range("O57").Value = "Región de Arica y Parinacota"
$ws.Range("P57").Value = 142
$ws.Range("Q57").Value = 60
$ws.Range("R57").Value = "Hortaliza"
